# chcked and fixed the ers
#
# The placements sheet had several duplicate company names in column B
# (Amazon x3, Accenture, IBM x2, Cognizant x2, Tech Mahindra) that made it
# hard to tell individual placement rows apart. Disambiguate the repeats by
# suffixing them with -2/-3, and add conditional formatting to flag any
# future duplicates in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Disambiguate duplicate company names -------------------------------
$ws.Range("B3").Value = "Amazon-2"
$ws.Range("B4").Value = "Amazon-3"
$ws.Range("B12").Value = "Accenture-2"
$ws.Range("B25").Value = "IBM-2"
$ws.Range("B28").Value = "IBM-3"
$ws.Range("B30").Value = "Cognizant-2"
$ws.Range("B31").Value = "Tech Mahindra-2"
$ws.Range("B32").Value = "Cognizant-3"

# --- Highlight duplicate values in column B ------------------------------
$col = $ws.Range("B1:B1048576")
$col.FormatConditions.Delete()
$fc = $col.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 0x06009C
$fc.Interior.Color = 0xCEC7FF

# --- Reset the selection back to the top of the sheet --------------------
$ws.Range("B1").Select()
